# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns with the
# latest scraped figures.
#
# Column D holds numeric-looking text (e.g. "1.001", "21.30", "0.00001032")
# that must stay text so values such as trailing zeros or the "." used as a
# thousands separator are preserved exactly as scraped. Excel auto-detects
# numbers when a value is assigned, so each target cell is switched to the
# "@" (Text) number format first to stop that conversion.
#
# Column E values (e.g. "  -3.80%  ") already contain spaces/percent signs
# so Excel keeps them as plain text without any extra formatting step.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.304.14"
$ws.Range("E2").Value = "  -3.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.851.10"
$ws.Range("E3").Value = "  -5.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.51%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.33"
$ws.Range("E5").Value = "  +0.72%  "

$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4500"
$ws.Range("E7").Value = "  -5.75%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3834"
$ws.Range("E8").Value = "  -4.77%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.27"
$ws.Range("E9").Value = "  -10.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07847"
$ws.Range("E10").Value = "  -7.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.014"
$ws.Range("E11").Value = "  -3.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.30"
$ws.Range("E12").Value = "  -4.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.848.09"
$ws.Range("E13").Value = "  -5.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.857"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.125"
$ws.Range("E15").Value = "  -5.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001032"
$ws.Range("E17").Value = "  -3.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "85.56"
$ws.Range("E18").Value = "  -5.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06509"
$ws.Range("E19").Value = "  -1.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.96"
$ws.Range("E20").Value = "  -8.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -0.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.473"
$ws.Range("E22").Value = "  -5.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.297.71"
$ws.Range("E23").Value = "  -3.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.80"
$ws.Range("E24").Value = "  -5.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.262"
$ws.Range("E25").Value = "  -0.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.055.42"
$ws.Range("E26").Value = "  -5.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.36"
$ws.Range("E27").Value = "  -2.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.37"
$ws.Range("E28").Value = "  -4.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.537"
$ws.Range("E29").Value = "  -6.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.049"
$ws.Range("E30").Value = "  -4.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.99"
$ws.Range("E31").Value = "  -3.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09313"
$ws.Range("E32").Value = "  -3.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.473"
$ws.Range("E33").Value = "  +1.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9331"
$ws.Range("E34").Value = "  -4.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.596"
$ws.Range("E35").Value = "  -1.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.267"
$ws.Range("E36").Value = "  -5.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02221"
$ws.Range("E37").Value = "  -4.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05972"
$ws.Range("E38").Value = "  -3.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.200"
$ws.Range("E39").Value = "  -4.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.229"
$ws.Range("E40").Value = "  -8.18%  "

$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5897"
$ws.Range("E42").Value = "  -4.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1854"
$ws.Range("E43").Value = "  -3.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.17"
$ws.Range("E44").Value = "  -8.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.253"
$ws.Range("E45").Value = "  -6.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5645"
$ws.Range("E46").Value = "  -5.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.10"
$ws.Range("E47").Value = "  -6.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.929"
$ws.Range("E48").Value = "  -6.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.357"
$ws.Range("E49").Value = "  -0.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06876"
$ws.Range("E50").Value = "  +1.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "108.01"
$ws.Range("E51").Value = "  -2.40%  "
